# Generate Report for Handback
#
# The localization status report is regenerated: the two source files
# (ac2f5080... and afb46ee5...) have been handed back from both the
# zh-cn and de-de locales, so their status flips from "Ready for
# handoff" to "Handed back: in sync with en-US", their "Latest Target
# File" / "Latest Handback File" columns get populated, and (for
# de-de, which just finished) "Latest Handback DateTime" gets a real
# timestamp.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$file1 = "ac2f5080-7f2a-49cf-be49-8ed8dacc307e.md"
$file2 = "afb46ee5-6896-4257-a56d-04be8f8c5f92.md"

$file1Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3da8bc095bc12110e89b16dcad1b50726d1bdebd/e2e/" + $file1
$file2Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3da8bc095bc12110e89b16dcad1b50726d1bdebd/e2e/" + $file2

# --- Overview sheet: status columns for both locales (E = zh-cn, F = de-de) ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $statusNew
$ov.Range("F2").Value = $statusNew
$ov.Range("E3").Value = $statusNew
$ov.Range("F3").Value = $statusNew

# Hyperlink-cell style: underlined text in the same blue used by the
# existing Source-File-Name hyperlinks (A2/A3) in this workbook.
$hyperlinkColor = 15570276  # BGR long for RGB FF6495ED ("cornflower blue")

function Set-HandbackHyperlink($sheetRange, $cellRef, $fileName, $fileUrl) {
    $cell = $sheetRange.Range($cellRef)
    $cell.Value = $fileName
    $sheetRange.Hyperlinks.Add($cell, $fileUrl, "", "", $fileName) | Out-Null
    $cell.Font.Name = "Calibri"
    $cell.Font.Underline = $true
    $cell.Font.Color = $hyperlinkColor
}

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $statusNew
$zh.Range("C3").Value = $statusNew

Set-HandbackHyperlink $zh "I2" $file1 $file1Url
$zh.Range("J2").Value = "ac2f5080-7f2a-49cf-be49-8ed8dacc307e.6271b15e4e1671a6ee414920087270d9c3b9af42.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-23 00:58:09"

Set-HandbackHyperlink $zh "I3" $file2 $file2Url
$zh.Range("J3").Value = "afb46ee5-6896-4257-a56d-04be8f8c5f92.03a0aa0ab94907ca03c25347035788a6465ddd12.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-23 00:58:09"

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $statusNew
$de.Range("C3").Value = $statusNew

Set-HandbackHyperlink $de "I2" $file1 $file1Url
$de.Range("J2").Value = "ac2f5080-7f2a-49cf-be49-8ed8dacc307e.6271b15e4e1671a6ee414920087270d9c3b9af42.de-de.xlf"
$de.Range("K2").Value = "2016-08-23 00:58:16"

Set-HandbackHyperlink $de "I3" $file2 $file2Url
$de.Range("J3").Value = "afb46ee5-6896-4257-a56d-04be8f8c5f92.03a0aa0ab94907ca03c25347035788a6465ddd12.de-de.xlf"
$de.Range("K3").Value = "2016-08-23 00:58:16"

# --- Column widths: widen the now-populated "Latest Target File" /
#     "Latest Handback File" columns (and the Status column, which
#     now holds a longer string) so the new values aren't clipped.
$ov.Columns.Item(5).ColumnWidth = 29.9777047293527
$ov.Columns.Item(6).ColumnWidth = 29.9777047293527

$zh.Columns.Item(3).ColumnWidth = 29.9777047293527
$zh.Columns.Item(9).ColumnWidth = 40
$zh.Columns.Item(10).ColumnWidth = 40

$de.Columns.Item(3).ColumnWidth = 29.9777047293527
$de.Columns.Item(9).ColumnWidth = 40
$de.Columns.Item(10).ColumnWidth = 40
